# Restore/update the rule table value in the "Rules" sheet.
# Cell C10 (row 10, the "R20" rule row) changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
